# Applies targeted numeric updates to the multiplication table cells
# by replacing each exact equation text with its new value.
$d = $word.ActiveDocument

$d.Content.Find.Execute("53×24=", $true, $false, $false, $false, $false, $true, 1, $false, "78×11=", 2) | Out-Null
$d.Content.Find.Execute("25×73=", $true, $false, $false, $false, $false, $true, 1, $false, "17×88=", 2) | Out-Null
$d.Content.Find.Execute("67×78=", $true, $false, $false, $false, $false, $true, 1, $false, "47×18=", 2) | Out-Null
$d.Content.Find.Execute("81×89=", $true, $false, $false, $false, $false, $true, 1, $false, "14×52=", 2) | Out-Null
$d.Content.Find.Execute("74×92=", $true, $false, $false, $false, $false, $true, 1, $false, "74×95=", 2) | Out-Null
$d.Content.Find.Execute("84×51=", $true, $false, $false, $false, $false, $true, 1, $false, "32×99=", 2) | Out-Null
$d.Content.Find.Execute("41×89=", $true, $false, $false, $false, $false, $true, 1, $false, "85×31=", 2) | Out-Null
$d.Content.Find.Execute("85×51=", $true, $false, $false, $false, $false, $true, 1, $false, "82×69=", 2) | Out-Null
$d.Content.Find.Execute("20×42=", $true, $false, $false, $false, $false, $true, 1, $false, "25×17=", 2) | Out-Null
$d.Content.Find.Execute("96×24=", $true, $false, $false, $false, $false, $true, 1, $false, "75×46=", 2) | Out-Null
$d.Content.Find.Execute("63×64=", $true, $false, $false, $false, $false, $true, 1, $false, "11×37=", 2) | Out-Null
$d.Content.Find.Execute("93×83=", $true, $false, $false, $false, $false, $true, 1, $false, "23×77=", 2) | Out-Null
$d.Content.Find.Execute("47×88=", $true, $false, $false, $false, $false, $true, 1, $false, "80×52=", 2) | Out-Null
$d.Content.Find.Execute("32×15=", $true, $false, $false, $false, $false, $true, 1, $false, "56×93=", 2) | Out-Null
$d.Content.Find.Execute("25×37=", $true, $false, $false, $false, $false, $true, 1, $false, "35×29=", 2) | Out-Null
$d.Content.Find.Execute("66×55=", $true, $false, $false, $false, $false, $true, 1, $false, "43×74=", 2) | Out-Null
$d.Content.Find.Execute("19×54=", $true, $false, $false, $false, $false, $true, 1, $false, "91×51=", 2) | Out-Null
$d.Content.Find.Execute("99×88=", $true, $false, $false, $false, $false, $true, 1, $false, "19×14=", 2) | Out-Null
$d.Content.Find.Execute("72×90=", $true, $false, $false, $false, $false, $true, 1, $false, "29×76=", 2) | Out-Null
$d.Content.Find.Execute("70×53=", $true, $false, $false, $false, $false, $true, 1, $false, "58×14=", 2) | Out-Null
$d.Content.Find.Execute("36×49=", $true, $false, $false, $false, $false, $true, 1, $false, "44×26=", 2) | Out-Null
$d.Content.Find.Execute("24×84=", $true, $false, $false, $false, $false, $true, 1, $false, "22×62=", 2) | Out-Null
$d.Content.Find.Execute("26×59=", $true, $false, $false, $false, $false, $true, 1, $false, "41×12=", 2) | Out-Null
$d.Content.Find.Execute("85×81=", $true, $false, $false, $false, $false, $true, 1, $false, "44×49=", 2) | Out-Null
$d.Content.Find.Execute("70×61=", $true, $false, $false, $false, $false, $true, 1, $false, "33×66=", 2) | Out-Null
